function Set-TextValue($ws, $ref, $val) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" '26.941.92'
Set-TextValue $ws "E2" '  +0.68%  '

Set-TextValue $ws "D3" '1.642.83'
Set-TextValue $ws "E3" '  -0.12%  '

Set-TextValue $ws "D4" '1.00'
Set-TextValue $ws "E4" '  -0.63%  '

Set-TextValue $ws "D5" '219.35'
Set-TextValue $ws "E5" '  +1.21%  '

Set-TextValue $ws "D6" '0.498'

Set-TextValue $ws "E7" '  -0.48%  '

Set-TextValue $ws "E8" '  -0.17%  '

Set-TextValue $ws "E9" '  -0.49%  '

Set-TextValue $ws "D10" '19.31'
Set-TextValue $ws "E10" '  +0.94%  '

Set-TextValue $ws "E11" '  +0.39%  '

Set-TextValue $ws "D12" '1.871.50'
Set-TextValue $ws "E12" '  -0.01%  '

Set-TextValue $ws "D13" '1.640.59'
Set-TextValue $ws "E13" '  -0.45%  '

Set-TextValue $ws "D14" '4.16'
Set-TextValue $ws "E14" '  -0.11%  '

Set-TextValue $ws "D15" '0.529'
Set-TextValue $ws "E15" '  +0.40%  '

Set-TextValue $ws "D16" '65.68'
Set-TextValue $ws "E16" '  +1.58%  '

Set-TextValue $ws "D17" '26.921.67'
Set-TextValue $ws "E17" '  +0.59%  '

Set-TextValue $ws "D18" '0.0₃0732'
Set-TextValue $ws "E18" '  -0.51%  '

Set-TextValue $ws "D19" '217.28'
Set-TextValue $ws "E19" '  +1.49%  '

Set-TextValue $ws "E20" '  -0.41%  '

Set-TextValue $ws "D21" '4.39'
Set-TextValue $ws "E21" '  -0.26%  '

Set-TextValue $ws "D22" '6.59'
Set-TextValue $ws "E22" '  +5.54%  '

Set-TextValue $ws "D23" '2.43'
Set-TextValue $ws "E23" '  -1.42%  '

Set-TextValue $ws "D24" '9.22'
Set-TextValue $ws "E24" '  -1.32%  '

Set-TextValue $ws "E25" '  +1.35%  '

Set-TextValue $ws "D26" '1.00'
Set-TextValue $ws "E26" '  -0.45%  '

Set-TextValue $ws "B27" 'Cosmos'
Set-TextValue $ws "C27" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws "D27" '7.30'
Set-TextValue $ws "E27" '  +2.40%  '

Set-TextValue $ws "B28" 'Stellar'
Set-TextValue $ws "C28" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws "D28" '0.119'
Set-TextValue $ws "E28" '  +0.28%  '

Set-TextValue $ws "D29" '15.83'
Set-TextValue $ws "E29" '  +1.32%  '

Set-TextValue $ws "D30" '0.0511'
Set-TextValue $ws "E30" '  +0.41%  '

Set-TextValue $ws "D31" '1.20'
Set-TextValue $ws "E31" '  +1.24%  '

Set-TextValue $ws "D32" '3.37'
Set-TextValue $ws "E32" '  -0.23%  '

Set-TextValue $ws "E33" '  -0.23%  '

Set-TextValue $ws "E34" '  +1.74%  '

Set-TextValue $ws "D35" '1.268.35'
Set-TextValue $ws "E35" '  -1.72%  '

Set-TextValue $ws "D36" '2.44'
Set-TextValue $ws "E36" '  -0.10%  '

Set-TextValue $ws "E37" '  -1.99%  '

Set-TextValue $ws "E38" '  -0.41%  '

Set-TextValue $ws "D39" '0.823'
Set-TextValue $ws "E39" '  +0.47%  '

Set-TextValue $ws "E40" '  -0.35%  '

Set-TextValue $ws "D41" '0.808'
Set-TextValue $ws "E41" '  +0.53%  '

Set-TextValue $ws "D42" '5.34'
Set-TextValue $ws "E42" '  +0.27%  '

Set-TextValue $ws "D43" '1.781.97'
Set-TextValue $ws "E43" '  -0.44%  '

Set-TextValue $ws "D44" '92.65'
Set-TextValue $ws "E44" '  +0.99%  '

Set-TextValue $ws "B45" 'Aave'
Set-TextValue $ws "C45" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws "D45" '61.05'
Set-TextValue $ws "E45" '  -1.03%  '

Set-TextValue $ws "B46" 'MXToken'
Set-TextValue $ws "C46" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws "D46" '2.08'
Set-TextValue $ws "E46" '  -6.44%  '

Set-TextValue $ws "E47" '  -0.34%  '

Set-TextValue $ws "D48" '0.0516'
Set-TextValue $ws "E48" '  -1.33%  '

Set-TextValue $ws "B49" 'BabyDogeCoin'
Set-TextValue $ws "C49" 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws "D49" '0.0₇0979'
Set-TextValue $ws "E49" '  -5.70%  '

Set-TextValue $ws "D50" '0.0971'
Set-TextValue $ws "E50" '  -0.01%  '

Set-TextValue $ws "B51" 'EnergySwap'
Set-TextValue $ws "C51" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws "D51" '7.59'
Set-TextValue $ws "E51" '  -1.00%  '
